# fixed plc parsing logic and wayside mapping in the excel
#
# Row 33's "Wayside" column (A33) incorrectly carried over the wayside1G
# value from the block above it. It should read wayside2G instead, matching
# the waysides used later in the sheet (rows 61+ already read wayside2G).
#
# Column A (rows 5-60) is one long chain of "=previous cell" formulas, so
# fixing the literal value in A33 automatically ripples the corrected
# wayside2G value down through A34:A60 (and the special A39 "=A35" cell)
# on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33").Value = "wayside2G"

# Restore the user's on-screen selection/scroll position to where they were
# working (near the bottom of the wayside2G block) when they made the fix.
$excel.ActiveWindow.ScrollRow = 112
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A124").Select()
